# C1--C2-and-C3-PowerPoint.pptx edit
#
# 1) Slide 16's table (graphic frame "Google Shape;213;p29", Shapes.Item(3))
#    switches from its table style {186513E6-EBD0-4FFB-A58A-CFE9C3158674}
#    (the deck's custom "Table_0" style) to the built-in gallery style
#    "Medium Style 2 - Accent 1" {4D6DE702-904B-4193-B7ED-6B63833AF73A}.
#
# 2) The presentation's theme color scheme (ppt/theme/theme1.xml, the
#    theme used by SlideMaster1 and therefore every slide) switches from
#    the custom "Integral" palette to the stock Office palette.

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{4D6DE702-904B-4193-B7ED-6B63833AF73A}")

# --- 2) Theme colors --------------------------------------------------------
# Office/Office Theme palette (replacing the "Integral" palette). RGB values
# are passed as COM OLE_COLOR longs (0xBBGGRR), matching PowerPoint's
# ThemeColorScheme.Item(n).RGB semantics. The theme is shared by the whole
# deck (SlideMaster1 / ppt/theme/theme1.xml), so any slide's
# ThemeColorScheme can be used to edit it - reuse the same slide as above.
$themeColors = $slide.ThemeColorScheme

$themeColors.Item(1).RGB  = 0         # dk1      000000
$themeColors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$themeColors.Item(3).RGB  = 6968388   # dk2      44546A
$themeColors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$themeColors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$themeColors.Item(6).RGB  = 3243501   # accent2  ED7D31
$themeColors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$themeColors.Item(8).RGB  = 49407     # accent4  FFC000
$themeColors.Item(9).RGB  = 12874308  # accent5  4472C4
$themeColors.Item(10).RGB = 4697456   # accent6  70AD47
$themeColors.Item(11).RGB = 12673797  # hlink    0563C1
$themeColors.Item(12).RGB = 7491477   # folHlink 954F72
